$d = $word.ActiveDocument

$replacements = @(
    @("920÷7=", "756÷2="),
    @("198÷5=", "452÷7="),
    @("118÷3=", "251÷4="),
    @("779÷3=", "330÷3="),
    @("829÷2=", "658÷8="),
    @("450÷8=", "341÷9="),
    @("362÷4=", "962÷3="),
    @("185÷9=", "558÷6="),
    @("401÷2=", "976÷5="),
    @("321÷6=", "751÷5="),
    @("304÷9=", "533÷9="),
    @("997÷4=", "840÷6="),
    @("686÷3=", "432÷2="),
    @("415÷3=", "674÷5="),
    @("590÷6=", "387÷9="),
    @("668÷4=", "781÷5="),
    @("339÷6=", "411÷3="),
    @("759÷2=", "552÷8="),
    @("302÷8=", "832÷6="),
    @("942÷3=", "295÷4="),
    @("656÷5=", "772÷3="),
    @("856÷2=", "422÷4="),
    @("525÷2=", "312÷4="),
    @("179÷8=", "618÷2="),
    @("171÷8=", "773÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
